$wb = $excel.ActiveWorkbook

$alc = $wb.Worksheets.Item("ALC")
$arm = $wb.Worksheets.Item("ARM")
$bsm = $wb.Worksheets.Item("BSM")
$crp = $wb.Worksheets.Item("CRP")
$cul = $wb.Worksheets.Item("CUL")
$gsm = $wb.Worksheets.Item("GSM")
$ltw = $wb.Worksheets.Item("LTW")
$wvr = $wb.Worksheets.Item("WVR")

# ALC row 41
$alc.Range("H41").Value = 797.6667
$alc.Range("I41").Value = 1008.8571
$alc.Range("J41").Value = 58.5
$alc.Range("K41").Value = 1008.8571
$alc.Range("L41").Value = 58.5
$alc.Range("M41").Value = -568.8570999999999
$alc.Range("N41").Value = -938.5

# ALC row 61
$alc.Range("H61").Value = 483.33334
$alc.Range("I61").Value = 483.33334
$alc.Range("K61").Value = 1450.00002
$alc.Range("M61").Value = -1278.00002

# ALC row 70
$alc.Range("H70").Value = 2699.2727
$alc.Range("I70").Value = 1599
$alc.Range("K70").Value = 4797
$alc.Range("M70").Value = -4527

# ALC row 73
$alc.Range("H73").Value = 2699.2727
$alc.Range("I73").Value = 1599
$alc.Range("K73").Value = 4797
$alc.Range("M73").Value = -3861

# ALC row 76
$alc.Range("H76").Value = 6081.1665
$alc.Range("I76").Value = 4871.75
$alc.Range("J76").Value = 8500
$alc.Range("K76").Value = 4871.75
$alc.Range("L76").Value = 8500
$alc.Range("M76").Value = -4556.75
$alc.Range("N76").Value = -9130

# ALC row 79
$alc.Range("H79").Value = 6081.1665
$alc.Range("I79").Value = 4871.75
$alc.Range("J79").Value = 8500
$alc.Range("K79").Value = 4871.75
$alc.Range("L79").Value = 8500
$alc.Range("M79").Value = -3779.75
$alc.Range("N79").Value = -10684

# ALC row 106
$alc.Range("H106").Value = 14224.5
$alc.Range("J106").Value = 15636
$alc.Range("L106").Value = 15636
$alc.Range("N106").Value = -16898

# ALC row 111
$alc.Range("H111").Value = 1495.4286
$alc.Range("I111").Value = 1475.6
$alc.Range("K111").Value = 4426.799999999999
$alc.Range("M111").Value = -1359.799999999999

# ALC row 118
$alc.Range("H118").Value = 2599.182
$alc.Range("I118").Value = 473.66666
$alc.Range("K118").Value = 1420.99998
$alc.Range("M118").Value = 236.0000199999999

# ALC row 132
$alc.Range("H132").Value = 8625.023999999999
$alc.Range("J132").Value = 5388
$alc.Range("L132").Value = 16164
$alc.Range("N132").Value = -21224

# ARM row 32
$arm.Range("H32").Value = 34588.855
$arm.Range("I32").Value = 38728.137
$arm.Range("K32").Value = 38728.137
$arm.Range("M32").Value = -38441.137

# ARM row 61
$arm.Range("H61").Value = 8783
$arm.Range("I61").Value = 7566.3335
$arm.Range("K61").Value = 7566.3335
$arm.Range("M61").Value = -7354.3335

# ARM row 63
$arm.Range("H63").Value = 2376.7144
$arm.Range("J63").Value = 1999.5
$arm.Range("L63").Value = 1999.5
$arm.Range("N63").Value = -3371.5

# ARM row 66
$arm.Range("H66").Value = 2376.7144
$arm.Range("J66").Value = 1999.5
$arm.Range("L66").Value = 9997.5
$arm.Range("N66").Value = -16861.5

# ARM row 122
$arm.Range("H122").Value = 1341.84
$arm.Range("I122").Value = 1294.9166
$arm.Range("K122").Value = 3884.7498
$arm.Range("M122").Value = -1434.7498

# ARM row 132
$arm.Range("H132").Value = 29656.395
$arm.Range("I132").Value = 33119.547
$arm.Range("J132").Value = 6799.6
$arm.Range("K132").Value = 99358.641
$arm.Range("L132").Value = 20398.8
$arm.Range("M132").Value = -96828.641
$arm.Range("N132").Value = -25458.8

# ARM row 136
$arm.Range("H136").Value = 8783
$arm.Range("I136").Value = 7566.3335
$arm.Range("K136").Value = 22699.0005
$arm.Range("M136").Value = -20149.0005

# BSM row 22
$bsm.Range("H22").Value = 125754.625
$bsm.Range("I22").Value = 125754.625
$bsm.Range("K22").Value = 125754.625
$bsm.Range("M22").Value = -125581.625

# BSM row 134
$bsm.Range("H134").Value = 2381.739
$bsm.Range("I134").Value = 2390.2222
$bsm.Range("K134").Value = 7170.6666
$bsm.Range("M134").Value = -4635.6666

# CRP row 16
$crp.Range("H16").Value = 1678
$crp.Range("I16").Value = 900
$crp.Range("J16").Value = 2283.111
$crp.Range("K16").Value = 900
$crp.Range("L16").Value = 2283.111
$crp.Range("M16").Value = -613
$crp.Range("N16").Value = -2857.111

# CRP row 58
$crp.Range("H58").Value = 336003.66
$crp.Range("I58").Value = 336003.66
$crp.Range("K58").Value = 336003.66
$crp.Range("M58").Value = -335800.66

# CRP row 113
$crp.Range("H113").Value = 1678
$crp.Range("I113").Value = 900
$crp.Range("J113").Value = 2283.111
$crp.Range("K113").Value = 900
$crp.Range("L113").Value = 2283.111
$crp.Range("M113").Value = 1270
$crp.Range("N113").Value = -6623.111

# CRP row 132
$crp.Range("H132").Value = 2732.484
$crp.Range("I132").Value = 2564.923
$crp.Range("K132").Value = 7694.768999999999
$crp.Range("M132").Value = -5164.768999999999

# CRP row 134
$crp.Range("H134").Value = 42270.84
$crp.Range("I134").Value = 54404.684
$crp.Range("K134").Value = 163214.052
$crp.Range("M134").Value = -160679.052

# CRP row 136
$crp.Range("H136").Value = 336003.66
$crp.Range("I136").Value = 336003.66
$crp.Range("K136").Value = 1008010.98
$crp.Range("M136").Value = -1005460.98

# CUL row 131
$cul.Range("H131").Value = 6270904
$cul.Range("J131").Value = 12528412
$cul.Range("L131").Value = 37585236
$cul.Range("N131").Value = -37595316

# CUL row 134
$cul.Range("H134").Value = 918.4
$cul.Range("I134").Value = 918.4
$cul.Range("K134").Value = 2755.2
$cul.Range("M134").Value = 2314.8

# CUL row 140
$cul.Range("H140").Value = 2956
$cul.Range("I140").Value = 2325.5
$cul.Range("K140").Value = 6976.5
$cul.Range("M140").Value = -1796.5

# GSM row 113
$gsm.Range("H113").Value = 145770.22
$gsm.Range("J113").Value = 253498.5
$gsm.Range("L113").Value = 253498.5
$gsm.Range("N113").Value = -257838.5

# GSM row 132
$gsm.Range("H132").Value = 503749
$gsm.Range("I132").Value = 999999
$gsm.Range("J132").Value = 7499
$gsm.Range("K132").Value = 2999997
$gsm.Range("L132").Value = 22497
$gsm.Range("M132").Value = -2997467
$gsm.Range("N132").Value = -27557

# LTW row 4
$ltw.Range("H4").Value = 0
$ltw.Range("J4").Value = 0
$ltw.Range("L4").Value = 0
$ltw.Range("N4").ClearContents()

# LTW row 7
$ltw.Range("H7").Value = 4148
$ltw.Range("I7").Value = 4148
$ltw.Range("K7").Value = 4148
$ltw.Range("M7").Value = -4036

# LTW row 16
$ltw.Range("H16").Value = 3235.0688
$ltw.Range("I16").Value = 3266.32
$ltw.Range("K16").Value = 3266.32
$ltw.Range("M16").Value = -3096.32

# LTW row 23
$ltw.Range("H23").Value = 17499
$ltw.Range("I23").Value = 17499
$ltw.Range("K23").Value = 17499
$ltw.Range("M23").Value = -17269

# LTW row 28
$ltw.Range("H28").Value = 0
$ltw.Range("J28").Value = 0
$ltw.Range("L28").Value = 0
$ltw.Range("N28").ClearContents()

# LTW row 37
$ltw.Range("H37").Value = 0
$ltw.Range("J37").Value = 0
$ltw.Range("L37").Value = 0
$ltw.Range("N37").ClearContents()

# LTW row 68
$ltw.Range("H68").Value = 6882.3335
$ltw.Range("I68").Value = 1950.5
$ltw.Range("J68").Value = 9348.25
$ltw.Range("K68").Value = 1950.5
$ltw.Range("L68").Value = 9348.25
$ltw.Range("M68").Value = -1201.5
$ltw.Range("N68").Value = -10846.25

# LTW row 71
$ltw.Range("H71").Value = 6882.3335
$ltw.Range("I71").Value = 1950.5
$ltw.Range("J71").Value = 9348.25
$ltw.Range("K71").Value = 9752.5
$ltw.Range("L71").Value = 46741.25
$ltw.Range("M71").Value = -6008.5
$ltw.Range("N71").Value = -54229.25

# LTW row 82
$ltw.Range("H82").Value = 2932.8333
$ltw.Range("J82").Value = 2866.6667
$ltw.Range("L82").Value = 2866.6667
$ltw.Range("N82").Value = -3588.6667

# LTW row 85
$ltw.Range("H85").Value = 2932.8333
$ltw.Range("J85").Value = 2866.6667
$ltw.Range("L85").Value = 2866.6667
$ltw.Range("N85").Value = -5362.6667

# LTW row 109
$ltw.Range("H109").Value = 69849.5
$ltw.Range("J109").Value = 69849.5
$ltw.Range("L109").Value = 69849.5
$ltw.Range("N109").Value = -72623.5

# LTW row 126
$ltw.Range("H126").Value = 4148
$ltw.Range("I126").Value = 4148
$ltw.Range("K126").Value = 12444
$ltw.Range("M126").Value = -9974

# LTW row 132
$ltw.Range("H132").Value = 36031.812
$ltw.Range("I132").Value = 46142.43
$ltw.Range("J132").Value = 4576.5557
$ltw.Range("K132").Value = 138427.29
$ltw.Range("L132").Value = 13729.6671
$ltw.Range("M132").Value = -135897.29
$ltw.Range("N132").Value = -18789.6671

# WVR row 62
$wvr.Range("H62").Value = 172415.5
$wvr.Range("I62").Value = 4331.6665
$wvr.Range("J62").Value = 340499.34
$wvr.Range("K62").Value = 4331.6665
$wvr.Range("L62").Value = 340499.34
$wvr.Range("M62").Value = -3707.6665
$wvr.Range("N62").Value = -341747.34

# WVR row 65
$wvr.Range("H65").Value = 172415.5
$wvr.Range("I65").Value = 4331.6665
$wvr.Range("J65").Value = 340499.34
$wvr.Range("K65").Value = 21658.3325
$wvr.Range("L65").Value = 1702496.7
$wvr.Range("M65").Value = -18538.3325
$wvr.Range("N65").Value = -1708736.7

# WVR row 81
$wvr.Range("H81").Value = 2076.3
$wvr.Range("I81").Value = 2016.75
$wvr.Range("K81").Value = 4033.5
$wvr.Range("M81").Value = -2972.5

# WVR row 84
$wvr.Range("H84").Value = 2076.3
$wvr.Range("I84").Value = 2016.75
$wvr.Range("K84").Value = 20167.5
$wvr.Range("M84").Value = -14863.5

# WVR row 132
$wvr.Range("H132").Value = 44602.707
$wvr.Range("I132").Value = 50688.906
$wvr.Range("J132").Value = 1999.3334
$wvr.Range("K132").Value = 152066.718
$wvr.Range("L132").Value = 5998.0002
$wvr.Range("M132").Value = -149536.718
$wvr.Range("N132").Value = -11058.0002
